$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# ------------------------------------------------------------------
# 1) Move the footer block (signature lines) from rows 22/23 down to
#    rows 24/25 to make room for the two extra worker/period rows
#    that are being inserted into the data table (16-19).
# ------------------------------------------------------------------

# Copy formatting (incl. number formats/borders/fill) onto the new rows.
# Only the populated B:C and H:J blocks exist on rows 22/23 (D:G are
# empty there), so copy those two blocks separately to avoid pulling in
# unrelated column-default formatting for D:G.
$ws.Range("B22:C22").Copy()
$ws.Range("B24:C24").PasteSpecial($xlPasteFormats)
$ws.Range("H22:J22").Copy()
$ws.Range("H24:J24").PasteSpecial($xlPasteFormats)
$ws.Range("B23:C23").Copy()
$ws.Range("B25:C25").PasteSpecial($xlPasteFormats)
$ws.Range("H23:J23").Copy()
$ws.Range("H25:J25").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# Re-create the merges on the new rows, and drop the old ones
$ws.Range("B22:C22").UnMerge()
$ws.Range("H22:J22").UnMerge()
$ws.Range("B23:C23").UnMerge()
$ws.Range("H23:J23").UnMerge()
$ws.Range("B24:C24").Merge()
$ws.Range("H24:J24").Merge()
$ws.Range("B25:C25").Merge()
$ws.Range("H25:J25").Merge()

# Carry over the text values
$ws.Range("B24").Value = "___________________________________"
$ws.Range("H24").Value = "___________________________________"
$ws.Range("B25").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H25").Value = "FIRMA DEL REPRESENTANTE LEGAL"

# Clear out the old rows 22/23 now that their content lives in 24/25
$ws.Range("B22:J23").Clear()

# ------------------------------------------------------------------
# 2) Expand the worker/period detail table (rows 16-17) into four
#    rows (16-19): each of the two workers now has two mora periods.
#    Row 16 keeps the "middle" banded style, row 19 keeps the
#    "bottom border" style that used to belong to row 17.
# ------------------------------------------------------------------

# Grab the bottom-border formatting (currently row 17) for the new
# last row (19) before it gets overwritten below.
$ws.Range("B17:J17").Copy()
$ws.Range("B19:J19").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# Fill rows 17 and 18 with the plain "middle" formatting taken from
# row 16.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial($xlPasteFormats)
$ws.Range("B18:J18").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# Row 16: ESTEFANY PEREZ CEBALLOS - period 2507
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1047477599"
$ws.Range("D16").Value = "ESTEFANY PEREZ CEBALLOS"
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 80000
$ws.Range("G16").Value = 2000000

# Row 17: ESTEFANY PEREZ CEBALLOS - period 2506
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1047477599"
$ws.Range("D17").Value = "ESTEFANY PEREZ CEBALLOS"
$ws.Range("E17").Value = "2506"
$ws.Range("F17").Value = 80000
$ws.Range("G17").Value = 2000000

# Row 18: LIZZETTE DEL ROSARIO HERRERA PEREZ - period 2507
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1047426490"
$ws.Range("D18").Value = "LIZZETTE DEL ROSARIO HERRERA PEREZ"
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# Row 19: LIZZETTE DEL ROSARIO HERRERA PEREZ - period 2506
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1047426490"
$ws.Range("D19").Value = "LIZZETTE DEL ROSARIO HERRERA PEREZ"
$ws.Range("E19").Value = "2506"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

# ------------------------------------------------------------------
# 3) Update the summary figures: total mora value and period count.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 273880
$ws.Range("F13").Value = 2

# ------------------------------------------------------------------
# 4) Let Excel recompute the "best fit" column widths now that the
#    data has changed.
# ------------------------------------------------------------------
$ws.Columns("B:J").AutoFit()
